# Auto-generated edit script: updates Price (D) / Volume(1h) (E) columns
# per the commit diff. Numeric-looking text values are forced to remain
# text (matching the source inlineStr data) without altering cell style,
# by briefly applying a Text number format and then resetting the style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.746.72"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.878.30"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4724"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3965"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08029"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "1.885.76"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.970"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.159"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001046"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06619"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "27.692.98"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.510"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "2.094.36"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.096"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.590"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9675"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09553"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.627"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.302"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.196"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5988"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1910"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5692"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.402"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06816"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("E51").Value = "  +8.46%  "
